$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage (no auto numeric conversion) for Price cells whose new value looks numeric,
# matching the source data which stores these as plain text strings.
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Apply the updated Price (D) and Volume(1h) (E) values
$ws.Range("D2").Value = '29.988.06'
$ws.Range("E2").Value = '  -0.13%  '
$ws.Range("D3").Value = '1.908.61'
$ws.Range("E3").Value = '  +0.30%  '
$ws.Range("D4").Value = '0.9997'
$ws.Range("E4").Value = '  -0.10%  '
$ws.Range("D5").Value = '0.7813'
$ws.Range("E5").Value = '  +4.81%  '
$ws.Range("D6").Value = '241.71'
$ws.Range("E7").Value = '  -0.03%  '
$ws.Range("D8").Value = '0.3156'
$ws.Range("E8").Value = '  +2.52%  '
$ws.Range("D9").Value = '26.15'
$ws.Range("E9").Value = '  +1.95%  '
$ws.Range("D10").Value = '0.06881'
$ws.Range("E10").Value = '  -0.48%  '
$ws.Range("D11").Value = '0.07967'
$ws.Range("E11").Value = '  -0.88%  '
$ws.Range("D12").Value = '1.910.63'
$ws.Range("E12").Value = '  +0.43%  '
$ws.Range("D13").Value = '0.7400'
$ws.Range("E13").Value = '  -2.24%  '
$ws.Range("E14").Value = '  -0.85%  '
$ws.Range("D15").Value = '92.65'
$ws.Range("E15").Value = '  +1.45%  '
$ws.Range("D16").Value = '29.998.64'
$ws.Range("E16").Value = '  -0.13%  '
$ws.Range("D17").Value = '13.92'
$ws.Range("E17").Value = '  -1.05%  '
$ws.Range("D18").Value = '5.866'
$ws.Range("E18").Value = '  -5.18%  '
$ws.Range("D19").Value = '245.27'
$ws.Range("E19").Value = '  +3.60%  '
$ws.Range("D20").Value = '0.000007729'
$ws.Range("E20").Value = '  -0.78%  '
$ws.Range("D21").Value = '0.9998'
$ws.Range("E21").Value = '  -0.05%  '
$ws.Range("D22").Value = '2.147.20'
$ws.Range("E22").Value = '  -0.43%  '
$ws.Range("D23").Value = '1.000'
$ws.Range("E23").Value = '  +0.00%  '
$ws.Range("D24").Value = '6.842'
$ws.Range("E24").Value = '  -3.56%  '
$ws.Range("D25").Value = '168.51'
$ws.Range("E25").Value = '  +0.42%  '
$ws.Range("D26").Value = '9.248'
$ws.Range("E26").Value = '  -1.05%  '
$ws.Range("D27").Value = '0.1376'
$ws.Range("E27").Value = '  +7.54%  '
$ws.Range("D28").Value = '18.85'
$ws.Range("D29").Value = '2.028'
$ws.Range("E29").Value = '  -1.15%  '
$ws.Range("D30").Value = '1.368'
$ws.Range("E30").Value = '  +1.32%  '
$ws.Range("E31").Value = '  -1.02%  '
$ws.Range("D32").Value = '4.293'
$ws.Range("E32").Value = '  -0.38%  '
$ws.Range("D33").Value = '4.071'
$ws.Range("E33").Value = '  +0.48%  '
$ws.Range("D34").Value = '0.05506'
$ws.Range("E34").Value = '  +4.08%  '
$ws.Range("D35").Value = '1.251'
$ws.Range("E35").Value = '  -2.68%  '
$ws.Range("D36").Value = '0.7321'
$ws.Range("E36").Value = '  -0.97%  '
$ws.Range("D37").Value = '2.728'
$ws.Range("E37").Value = '  +0.06%  '
$ws.Range("D38").Value = '0.01927'
$ws.Range("E38").Value = '  -1.13%  '
$ws.Range("D39").Value = '2.788'
$ws.Range("E39").Value = '  +0.65%  '
$ws.Range("D40").Value = '6.124'
$ws.Range("E40").Value = '  -2.12%  '
$ws.Range("D41").Value = '0.4410'
$ws.Range("D42").Value = '71.82'
$ws.Range("E42").Value = '  -1.05%  '
$ws.Range("D43").Value = '1.000'
$ws.Range("E43").Value = '  +0.02%  '
$ws.Range("D44").Value = '0.8404'
$ws.Range("E44").Value = '  +0.98%  '
$ws.Range("D45").Value = '1.866'
$ws.Range("E45").Value = '  -4.37%  '
$ws.Range("D46").Value = '100.34'
$ws.Range("E46").Value = '  -1.12%  '
$ws.Range("D47").Value = '7.508'
$ws.Range("E47").Value = '  -2.94%  '
$ws.Range("E48").Value = '  -0.84%  '
$ws.Range("D49").Value = '984.57'
$ws.Range("E49").Value = '  +8.44%  '
$ws.Range("D50").Value = '2.056.51'
$ws.Range("E50").Value = '  +0.04%  '
$ws.Range("D51").Value = '36.20'
$ws.Range("E51").Value = '  -1.04%  '
